# Insert a new data row before the current row 42, shifting all the
# existing rows (42..122) down to (43..123). Excel copies the formatting
# of the row above into the newly inserted row, which matches the source
# workbook (the D column keeps its date number format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(42).Insert()

# Populate the freshly inserted row 42 with the new market-price record.
$ws.Cells.Item(42, 1).Value  = 1
$ws.Cells.Item(42, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(42, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(42, 4).Value  = 44952
$ws.Cells.Item(42, 5).Value  = 15
$ws.Cells.Item(42, 6).Value  = "Fruta"
$ws.Cells.Item(42, 7).Value  = 100102
$ws.Cells.Item(42, 8).Value  = "Cítricos"
$ws.Cells.Item(42, 9).Value  = 100102005
$ws.Cells.Item(42, 10).Value = "Naranja"
$ws.Cells.Item(42, 11).Value = "Valencia"
$ws.Cells.Item(42, 12).Value = "Segunda"
$ws.Cells.Item(42, 13).Value = 350
$ws.Cells.Item(42, 14).Value = 2200
$ws.Cells.Item(42, 15).Value = 2400
$ws.Cells.Item(42, 16).Value = 2286
$ws.Cells.Item(42, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(42, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(42, 19).Value = 2286
$ws.Cells.Item(42, 20).Value = 1
